$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I and J, rows 2-11
$values = @(
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 2),
    @(1, 6),
    @(1, 6),
    @(6, 6),
    @(1, 5),
    @(1, 5),
    @(8, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
